# Daily scrape update - 2025-08-07 03:43:45 UTC
# Replaces rows 2-8 with freshly scraped opportunities, drops the old
# rows 9-10 (list shrank from 10 to 8 entries), and adjusts a few
# column widths to fit the new content.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value into column A (the opportunity id) as TEXT -
# the source data is a numeric-looking string and must stay text (it
# is stored as inlineStr in the workbook, not a number). Tagging the
# cell "@" (Text) before assigning keeps the digits from being parsed
# as a number, and resetting the style back to Normal afterwards keeps
# the cell formatting identical to the surrounding, untouched cells.
function Set-IdText($addr, $val) {
    $c = $ws.Range($addr)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.Style = "Normal"
}

# --- Row 2 ---------------------------------------------------------
Set-IdText "A2" "1326671"
$ws.Range("B2").Value = "https://aiesec.org/opportunity/global-talent/1326671"
$ws.Range("C2").Value = "Marketing Intern"
$ws.Range("D2").Value = "Mumbai, Maharashtra, India"
$ws.Range("E2").Value = "No"
$ws.Range("F2").Value = "0 applicants"
$ws.Range("G2").Value = "3 - 6 Months"
$ws.Range("H2").Value = "India Someday Travels LLP"

# --- Row 3 ---------------------------------------------------------
Set-IdText "A3" "1326475"
$ws.Range("B3").Value = "https://aiesec.org/opportunity/global-talent/1326475"
$ws.Range("C3").Value = "Business development intern"
$ws.Range("D3").Value = "Chandigarh, India"
$ws.Range("E3").Value = "No"
$ws.Range("F3").Value = "1 applicant"
$ws.Range("G3").Value = "3 - 6 Months"
$ws.Range("H3").Value = "AgNext Technologies Private ltd"

# --- Row 4 ---------------------------------------------------------
Set-IdText "A4" "1326337"
$ws.Range("B4").Value = "https://aiesec.org/opportunity/global-talent/1326337"
$ws.Range("C4").Value = "Digital Marketing"
$ws.Range("D4").Value = "Lisboa, Portugal"
$ws.Range("E4").Value = "No"
$ws.Range("F4").Value = "34 applicants"
$ws.Range("G4").Value = "6 - 18 Months"
$ws.Range("H4").Value = "BGI S.A"

# --- Row 5 ---------------------------------------------------------
Set-IdText "A5" "1325876"
$ws.Range("B5").Value = "https://aiesec.org/opportunity/global-talent/1325876"
$ws.Range("C5").Value = "[Impact Brazil]- AI Data Scientist Intern"
$ws.Range("D5").Value = "Ribeirão Preto, SP, Brasil"
$ws.Range("E5").Value = "No"
$ws.Range("F5").Value = "71 applicants"
$ws.Range("G5").Value = "6 - 18 Months"
$ws.Range("H5").Value = "CCM Soluções em Tecnologia LTDA"

# --- Row 6 ---------------------------------------------------------
Set-IdText "A6" "1324636"
$ws.Range("B6").Value = "https://aiesec.org/opportunity/global-talent/1324636"
$ws.Range("C6").Value = "[Impact Fortaleza] -Cost & Quality Planning"
$ws.Range("D6").Value = "Castanhal, Pará, Brasil"
$ws.Range("E6").Value = "No"
$ws.Range("F6").Value = "8 applicants"
$ws.Range("G6").Value = "6 - 18 Months"
$ws.Range("H6").Value = "Petruz Fruity"

# --- Row 7 ---------------------------------------------------------
Set-IdText "A7" "1322455"
$ws.Range("B7").Value = "https://aiesec.org/opportunity/global-talent/1322455"
$ws.Range("C7").Value = "[Impact Fortaleza]- Chemical Engineering"
$ws.Range("D7").Value = "Castanhal, PA, Brasil"
$ws.Range("E7").Value = "No"
$ws.Range("F7").Value = "16 applicants"
$ws.Range("G7").Value = "6 - 18 Months"
$ws.Range("H7").Value = "Petruz Fruity"

# --- Row 8 ---------------------------------------------------------
Set-IdText "A8" "1315961"
$ws.Range("B8").Value = "https://aiesec.org/opportunity/global-talent/1315961"
$ws.Range("C8").Value = "Female Guest Relations Executive - Mid Term"
$ws.Range("D8").Value = "Kandy, Sri Lanka"
$ws.Range("E8").Value = "No"
$ws.Range("F8").Value = "23 applicants"
$ws.Range("G8").Value = "3 - 6 Months"
$ws.Range("H8").Value = "Canora Hotels (pvt) Ltd Grand Kandyan Hotel"

# --- Rows 9-10 no longer present in the scrape: drop them -----------
$ws.Range("A9:H10").EntireRow.Delete()

# --- Column width tweaks to fit the new content ---------------------
$ws.Columns.Item(3).ColumnWidth = 45.166666666666664   # C: 61 -> 46
$ws.Columns.Item(4).ColumnWidth = 28.166666666666668   # D: 78 -> 29
$ws.Columns.Item(6).ColumnWidth = 15.166666666666666   # F: 17 -> 16
$ws.Columns.Item(8).ColumnWidth = 45.166666666666664   # H: 29 -> 46
